$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.566.31"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.238.02"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'604.85"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "'156.91"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.238.25"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("D11").Value = "'5.79"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").Value = "'39.06"
$ws.Range("E14").Value = "  +0.42%  "
$ws.Range("D15").Value = "3.771.38"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "66.603.08"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "3.246.64"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'508.58"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("D21").Value = "'15.29"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "'0.746"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'8.03"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'14.67"
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("D25").Value = "'86.03"
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "'0.169"
$ws.Range("E26").Value = "  +88.77%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").Value = "'9.08"
$ws.Range("E29").Value = "  -1.21%  "
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").Value = "'6.91"
$ws.Range("E32").Value = "  -1.02%  "
$ws.Range("D33").Value = "'28.16"
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("E35").Value = "  -4.45%  "
$ws.Range("D36").Value = "'6.37"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "0.0₃0807"
$ws.Range("E37").Value = "  +18.98%  "
$ws.Range("D38").Value = "'55.30"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "'495.46"
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("E40").Value = "  +13.71%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").Value = "'8.75"
$ws.Range("E43").Value = "  -1.19%  "
$ws.Range("E44").Value = "  -3.50%  "
$ws.Range("D45").Value = "2.948.06"
$ws.Range("E45").Value = "  +3.23%  "
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").Value = "'28.25"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -1.07%  "
